$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("By year")
$ws.Range("G2").Value = 71.79
$ws.Range("H2").Value = 501.5
$ws.Range("I2").Value = 953.7
$ws.Range("J2").Value = 0.1636
$ws.Range("K2").Value = 866.8
$ws.Range("L2").Value = 1224
$ws.Range("M2").Value = 5216
$ws.Range("N2").Value = 6543
$ws.Range("O2").Value = 0.1503
$ws.Range("P2").Value = 0.2047
$ws.Range("G3").Value = 81.44
$ws.Range("H3").Value = 345.6
$ws.Range("I3").Value = 810
$ws.Range("J3").Value = 0.1371
$ws.Range("K3").Value = 668.9
$ws.Range("L3").Value = 1040
$ws.Range("M3").Value = 5287
$ws.Range("N3").Value = 6632
$ws.Range("O3").Value = 0.1134
$ws.Range("P3").Value = 0.164
$ws.Range("G4").Value = 71.79
$ws.Range("H4").Value = 369
$ws.Range("I4").Value = 697
$ws.Range("J4").Value = 0.1196
$ws.Range("K4").Value = 529.4
$ws.Range("L4").Value = 881.2
$ws.Range("M4").Value = 5216
$ws.Range("N4").Value = 6543
$ws.Range("O4").Value = 0.09052
$ws.Range("P4").Value = 0.1406

$ws = $wb.Worksheets.Item("By year, livestock cat.")
$ws.Range("H2").Value = 71.79
$ws.Range("I2").Value = 249.9
$ws.Range("J2").Value = 321.7
$ws.Range("K2").Value = 0.3007
$ws.Range("L2").Value = 223.4
$ws.Range("M2").Value = 417
$ws.Range("N2").Value = 957.3
$ws.Range("O2").Value = 1201
$ws.Range("P2").Value = 0.2093
$ws.Range("Q2").Value = 0.3822
$ws.Range("H3").Value = 130.4
$ws.Range("I3").Value = 501.5
$ws.Range("J3").Value = 632
$ws.Range("K3").Value = 0.1328
$ws.Range("L3").Value = 592.5
$ws.Range("M3").Value = 867.1
$ws.Range("N3").Value = 4258
$ws.Range("O3").Value = 5342
$ws.Range("P3").Value = 0.1253
$ws.Range("Q3").Value = 0.1746
$ws.Range("H4").Value = 81.44
$ws.Range("I4").Value = 260.1
$ws.Range("J4").Value = 341.6
$ws.Range("K4").Value = 0.2406
$ws.Range("L4").Value = 258.8
$ws.Range("M4").Value = 457.8
$ws.Range("N4").Value = 1270
$ws.Range("O4").Value = 1594
$ws.Range("P4").Value = 0.1761
$ws.Range("Q4").Value = 0.3133
$ws.Range("H5").Value = 122.8
$ws.Range("I5").Value = 345.6
$ws.Range("J5").Value = 468.5
$ws.Range("K5").Value = 0.1043
$ws.Range("L5").Value = 367.2
$ws.Range("M5").Value = 599.2
$ws.Range("N5").Value = 4017
$ws.Range("O5").Value = 5039
$ws.Range("P5").Value = 0.083
$ws.Range("Q5").Value = 0.1259
$ws.Range("H6").Value = 71.79
$ws.Range("I6").Value = 141.2
$ws.Range("J6").Value = 213
$ws.Range("K6").Value = 0.1991
$ws.Range("L6").Value = 174.1
$ws.Range("M6").Value = 253.6
$ws.Range("N6").Value = 957.3
$ws.Range("O6").Value = 1201
$ws.Range("P6").Value = 0.1576
$ws.Range("Q6").Value = 0.2305
$ws.Range("H7").Value = 115
$ws.Range("I7").Value = 369
$ws.Range("J7").Value = 484
$ws.Range("K7").Value = 0.1017
$ws.Range("L7").Value = 350.3
$ws.Range("M7").Value = 616.3
$ws.Range("N7").Value = 4258
$ws.Range("O7").Value = 5342
$ws.Range("P7").Value = 0.07356
$ws.Range("Q7").Value = 0.1253

$ws = $wb.Worksheets.Item("By application")
$ws.Range("AB2").Value = 18.46
$ws.Range("AC2").Value = 0.2051
$ws.Range("AB3").Value = 31.24
$ws.Range("AC3").Value = 0.3471
$ws.Range("AB4").Value = 12.62
$ws.Range("AC4").Value = 0.1052
$ws.Range("AB5").Value = 17.1
$ws.Range("AC5").Value = 0.1425
$ws.Range("AB6").Value = 20.36
$ws.Range("AC6").Value = 0.1939
$ws.Range("AB7").Value = 27.32
$ws.Range("AC7").Value = 0.2601
$ws.Range("AB8").Value = 13.94
$ws.Range("AC8").Value = 0.1032
$ws.Range("AB9").Value = 12.57
$ws.Range("AC9").Value = 0.1047
$ws.Range("AB10").Value = 18.46
$ws.Range("AC10").Value = 0.2051
$ws.Range("AB11").Value = 17.65
$ws.Range("AC11").Value = 0.1961
$ws.Range("AB12").Value = 11.13
$ws.Range("AC12").Value = 0.09276
$ws.Range("AB13").Value = 12.58
$ws.Range("AC13").Value = 0.1048
